# Progress update as of 04-Nov-2025:
# - "PERIOD TO EXPIRE" (col H) drops by one day for every training row.
# - "LAST UPDATE" (col I) moves from 03-Nov-2025 to 04-Nov-2025.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 20; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    $hCell.Value = $hCell.Value2 - 1

    # Leading apostrophe keeps this a literal text date (matching the
    # existing "03-Nov-2025" text entries) instead of Excel's automatic
    # text-to-date-serial conversion.
    $iCell.Value = "'04-Nov-2025"
}
